$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 14:52"

# Row 41 - Arabia Saudita: refresh case counts (column A/D/F unchanged)
$ws.Range("B41").Value = 2932
$ws.Range("C41").Value = 137
$ws.Range("E41").Value = 2276

# Rows 55/56 - Singapur overtakes Islandia in the sorted ranking.
# Row 55 becomes Singapur with its refreshed figures; row 56 becomes
# Islandia, keeping the figures it already had.
$ws.Range("A55").Value = "Singapur"
$ws.Range("B55").Value = 1623
$ws.Range("C55").Value = 142
$ws.Range("D55").Value = 406
$ws.Range("E55").Value = 1211
$ws.Range("F55").Value = 29

$ws.Range("A56").Value = "Islandia"
$ws.Range("B56").Value = 1586
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 559
$ws.Range("E56").Value = 1021
$ws.Range("F56").Value = 11

# Row 72 - Barein: refresh case counts (column A/D unchanged value, F unchanged)
$ws.Range("B72").Value = 821
$ws.Range("C72").Value = 10
$ws.Range("D72").Value = 467
$ws.Range("E72").Value = 349

# Row 116 - Kenia: refresh case counts (F unchanged)
$ws.Range("B116").Value = 179
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 9
$ws.Range("E116").Value = 164
